# Apply cell-value updates for the cryptos list refresh.
# Each target cell is forced to Text format before/while assigning the
# value so that numeric-looking strings (e.g. "1.00", "0.0000220") are
# preserved verbatim instead of being normalized/rounded by Excel, then
# the style is reset back to "Normal" so no stray number-format style is
# left attached to the cell (matching the original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '62.327.66'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.58%  '
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.998.82'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -0.90%  '
$c.Style = "Normal"
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '543.54'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -1.87%  '
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '138.00'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +1.60%  '
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '2.994.29'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.488'
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -1.73%  '
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +12.73%  '
$c.Style = "Normal"
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.148'
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -1.08%  '
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.446'
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.97%  '
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.0000220'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c.Style = "Normal"
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '33.96'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -1.17%  '
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.466.17'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '62.345.08'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.41%  '
$c.Style = "Normal"
$c = $ws.Range('B17')
$c.NumberFormat = "@"
$c.Value = 'WrappedEther'
$c.Style = "Normal"
$c = $ws.Range('C17')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.996.34'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -1.08%  '
$c.Style = "Normal"
$c = $ws.Range('B18')
$c.NumberFormat = "@"
$c.Value = 'TRON'
$c.Style = "Normal"
$c = $ws.Range('C18')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -2.17%  '
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '6.57'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -1.61%  '
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '469.49'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -0.91%  '
$c.Style = "Normal"
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.40'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.654'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -3.27%  '
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '7.16'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.99%  '
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '79.31'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -1.38%  '
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '12.57'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +3.47%  '
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.72'
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -0.70%  '
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +4.66%  '
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '25.43'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -1.46%  '
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -2.51%  '
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '2.34'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +0.40%  '
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.54'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +1.88%  '
$c.Style = "Normal"
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '54.67'
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -1.82%  '
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -1.49%  '
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '451.26'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -2.09%  '
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0810'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +1.54%  '
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.0392'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +1.74%  '
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '2.947.76'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -8.35%  '
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -3.78%  '
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '8.04'
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -1.15%  '
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +2.74%  '
$c.Style = "Normal"
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '26.77'
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +3.60%  '
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.248'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +1.17%  '
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +0.33%  '
$c.Style = "Normal"
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.99'
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -0.14%  '
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '114.98'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -2.80%  '
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0₃0494'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c.Style = "Normal"
$c = $ws.Range('B51')
$c.NumberFormat = "@"
$c.Value = 'ThetaToken'
$c.Style = "Normal"
$c = $ws.Range('C51')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.01'
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.66%  '
$c.Style = "Normal"
